$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 503.20514
$ws.Range("I15").Value = 503.20514
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 1509.61542
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -1340.61542

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 1542.8
$ws.Range("I99").Value = 1542.8
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4628.4
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -3130.4
$ws.Range("N99").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2284.2632
$ws.Range("I100").Value = 2284.2632
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2284.2632
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1743.2632

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 3357.5
$ws.Range("I103").Value = 3000
$ws.Range("J103").Value = 4668.3335
$ws.Range("K103").Value = 9000
$ws.Range("L103").Value = 14005.0005
$ws.Range("M103").Value = -8414
$ws.Range("N103").Value = -15177.0005

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 48500
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 48500
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 48500
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -58300

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5641.65
$ws.Range("I138").Value = 4694.5454
$ws.Range("J138").Value = 6799.222
$ws.Range("K138").Value = 14083.6362
$ws.Range("L138").Value = 20397.666
$ws.Range("M138").Value = -8943.636200000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3254.5557
$ws.Range("I132").Value = 2870.1428
$ws.Range("J132").Value = 4600
$ws.Range("K132").Value = 8610.428400000001
$ws.Range("L132").Value = 13800
$ws.Range("M132").Value = -6080.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 148.5
$ws.Range("I64").Value = 97.5
$ws.Range("J64").Value = 199.5
$ws.Range("K64").Value = 97.5
$ws.Range("L64").Value = 199.5
$ws.Range("M64").Value = 127.5
$ws.Range("N64").Value = -649.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 148.5
$ws.Range("I67").Value = 97.5
$ws.Range("J67").Value = 199.5
$ws.Range("K67").Value = 97.5
$ws.Range("L67").Value = 199.5
$ws.Range("M67").Value = 682.5
$ws.Range("N67").Value = -1759.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 473.1
$ws.Range("I80").Value = 432.8
$ws.Range("J80").Value = 513.4
$ws.Range("K80").Value = 432.8
$ws.Range("L80").Value = 513.4
$ws.Range("M80").Value = 565.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 473.1
$ws.Range("I83").Value = 432.8
$ws.Range("J83").Value = 513.4
$ws.Range("K83").Value = 2164
$ws.Range("L83").Value = 2567
$ws.Range("M83").Value = 2828

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3624.0833
$ws.Range("I105").Value = 3248.9
$ws.Range("J105").Value = 5500
$ws.Range("K105").Value = 3248.9
$ws.Range("L105").Value = 5500
$ws.Range("M105").Value = -1501.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3368.25
$ws.Range("I134").Value = 1739.5
$ws.Range("J134").Value = 4997
$ws.Range("K134").Value = 5218.5
$ws.Range("L134").Value = 14991
$ws.Range("M134").Value = -2683.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 565.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 896.1667
$ws.Range("I16").Value = 895.6
$ws.Range("J16").Value = 899
$ws.Range("K16").Value = 895.6
$ws.Range("L16").Value = 899
$ws.Range("M16").Value = -608.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 408
$ws.Range("I25").Value = 310
$ws.Range("J25").Value = 800
$ws.Range("K25").Value = 310
$ws.Range("L25").Value = 800
$ws.Range("M25").Value = -136

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5861
$ws.Range("I58").Value = 5035.3335
$ws.Range("J58").Value = 8338
$ws.Range("K58").Value = 5035.3335
$ws.Range("L58").Value = 8338
$ws.Range("M58").Value = -4832.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 6538.222
$ws.Range("I99").Value = 5637.8335
$ws.Range("J99").Value = 8339
$ws.Range("K99").Value = 5637.8335
$ws.Range("L99").Value = 8339
$ws.Range("M99").Value = -4139.8335
$ws.Range("N99").Value = -11335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1999.2727
$ws.Range("I105").Value = 1498.6
$ws.Range("J105").Value = 2416.5
$ws.Range("K105").Value = 1498.6
$ws.Range("L105").Value = 2416.5
$ws.Range("M105").Value = 248.4000000000001
$ws.Range("N105").Value = -5910.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 896.1667
$ws.Range("I113").Value = 895.6
$ws.Range("J113").Value = 899
$ws.Range("K113").Value = 895.6
$ws.Range("L113").Value = 899
$ws.Range("M113").Value = 1274.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 6538.222
$ws.Range("I126").Value = 5637.8335
$ws.Range("J126").Value = 8339
$ws.Range("K126").Value = 16913.5005
$ws.Range("L126").Value = 25017
$ws.Range("M126").Value = -14443.5005
$ws.Range("N126").Value = -29957

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2045.4117
$ws.Range("I134").Value = 1561.2858
$ws.Range("J134").Value = 4304.6665
$ws.Range("K134").Value = 4683.857400000001
$ws.Range("L134").Value = 12913.9995
$ws.Range("M134").Value = -2148.857400000001
$ws.Range("N134").Value = -17983.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 5861
$ws.Range("I136").Value = 5035.3335
$ws.Range("J136").Value = 8338
$ws.Range("K136").Value = 15106.0005
$ws.Range("L136").Value = 25014
$ws.Range("M136").Value = -12556.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 10045.333
$ws.Range("I11").Value = 138
$ws.Range("J11").Value = 14999
$ws.Range("K11").Value = 414
$ws.Range("L11").Value = 44997
$ws.Range("M11").Value = -274

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 95
$ws.Range("I13").Value = 95
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 285
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -117

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 891.2857
$ws.Range("I23").Value = 856.3333
$ws.Range("J23").Value = 917.5
$ws.Range("K23").Value = 2568.9999
$ws.Range("L23").Value = 2752.5
$ws.Range("M23").Value = -2333.9999
$ws.Range("N23").Value = -3222.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 1499.6666
$ws.Range("I41").Value = 1499.5
$ws.Range("J41").Value = 1500
$ws.Range("K41").Value = 4498.5
$ws.Range("L41").Value = 4500
$ws.Range("M41").Value = -4160.5
$ws.Range("N41").Value = -5176

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 900
$ws.Range("I75").Value = 900
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 2700
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -1702

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 900
$ws.Range("I78").Value = 900
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 8100
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -3108

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 841.25
$ws.Range("I86").Value = 727.8570999999999
$ws.Range("J86").Value = 1000
$ws.Range("K86").Value = 2183.5713
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -997.5712999999996
$ws.Range("N86").Value = -5372

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 841.25
$ws.Range("I89").Value = 727.8570999999999
$ws.Range("J89").Value = 1000
$ws.Range("K89").Value = 6550.7139
$ws.Range("L89").Value = 9000
$ws.Range("M89").Value = -622.7138999999997
$ws.Range("N89").Value = -20856

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1510.65
$ws.Range("I107").Value = 599.8
$ws.Range("J107").Value = 1814.2667
$ws.Range("K107").Value = 1799.4
$ws.Range("L107").Value = 5442.800099999999
$ws.Range("M107").Value = 120.6000000000001
$ws.Range("N107").Value = -9282.8001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3550

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4615.25
$ws.Range("I132").Value = 4837.3335
$ws.Range("J132").Value = 3949
$ws.Range("K132").Value = 14512.0005
$ws.Range("L132").Value = 11847
$ws.Range("M132").Value = -11982.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 6336636.5
$ws.Range("I100").Value = 8712500
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 17425000
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -17424459

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2001.3334
$ws.Range("I126").Value = 1999.5
$ws.Range("J126").Value = 2005
$ws.Range("K126").Value = 5998.5
$ws.Range("L126").Value = 6015
$ws.Range("M126").Value = -3528.5
$ws.Range("N126").Value = -10955

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1855.9166
$ws.Range("I132").Value = 1855.9166
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5567.7498
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3037.7498
$ws.Range("N132").ClearContents()
